$d = $word.ActiveDocument

# The document body currently ends with a paragraph that holds the big
# "Class Diagram" drawing (an AlternateContent picture), immediately
# followed by the section properties (sectPr). We need to append, right
# after that drawing paragraph:
#   1. a new, empty paragraph
#   2. a new Heading1 paragraph with the text "Controller State Diagram"
#
# Note: Word's Paragraphs collection also enumerates the paragraphs
# nested inside the drawing's text boxes; those all report a
# zero-length range collapsed onto the end of the document, so we
# can't just use Paragraphs.Count/Last after inserting. Instead we
# find the real trailing (non-collapsed) paragraph first and then
# address the newly inserted paragraphs by a fixed offset from it.

$count = $d.Paragraphs.Count
$drawingIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -ne $candidate.Range.End) {
        $drawingIndex = $i
        break
    }
}

$drawingPara = $d.Paragraphs.Item($drawingIndex)
$drawingPara.Range.InsertParagraphAfter()

$emptyPara = $d.Paragraphs.Item($drawingIndex + 1)
$emptyPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($drawingIndex + 2)
$headingPara.Range.Text = "Controller State Diagram"
$headingPara.Style = "Heading 1"
